# Actualización automática 2025-08-05 17:15:08
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("H3").Value = 267.3
$wsGrupo.Range("M54").Value = 245.95
$wsGrupo.Range("O54").Value = 637.42
$wsGrupo.Range("H56").Value = "1 de 54"
$wsGrupo.Range("O56").Value = "1 de 54"

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F3").Value = 267.3
$wsMensual.Range("F54").Value = 883.37
$wsMensual.Range("F56").Value = 13292.56

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D7").Value = 267.3
$wsCumpl.Range("E7").Value = 2132.7
$wsCumpl.Range("F7").Value = 0.111375

$wsCumpl.Range("D16").Value = 10820.89
$wsCumpl.Range("E16").Value = 45238.81
$wsCumpl.Range("F16").Value = 0.1930244007727476

$wsCumpl.Range("D18").Value = 637.42
$wsCumpl.Range("E18").Value = 2562.58
$wsCumpl.Range("F18").Value = 0.19919375

$wsCumpl.Range("D19").Value = 13292.56
$wsCumpl.Range("E19").Value = 104147.1306451791
$wsCumpl.Range("F19").Value = 0.1131862654522894
